# Update the "Förändrad" (Changed) date column (C) for every data row
# (rows 2-418) from serial date 46075 to 46076 (bump the "last changed"
# date by one day), matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 418; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46075) {
        $cell.Value = 46076
    }
}
